$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.490.94'
$ws.Range('E2').Value = '  +6.43%  '
$ws.Range('D3').Value = '3.322.69'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '409.98'
$ws.Range('E5').Value = '  +3.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.39'
$ws.Range('E6').Value = '  +5.94%  '
$ws.Range('D7').Value = '3.317.49'
$ws.Range('E7').Value = '  +2.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -1.61%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.626'
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('E11').Value = '  +17.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '39.98'
$ws.Range('E12').Value = '  +1.88%  '
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '3.847.02'
$ws.Range('E14').Value = '  +2.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.20'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.13'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '3.323.22'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').Value = '60.414.24'
$ws.Range('E18').Value = '  +6.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.00'
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.79'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000114'
$ws.Range('E21').Value = '  +4.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.37'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('E23').Value = '  -3.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '295.17'
$ws.Range('E24').Value = '  -0.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.89'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.64'
$ws.Range('E28').Value = '  +4.38%  '
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.50'
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('E32').Value = '  +4.47%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.27'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.09'
$ws.Range('E35').Value = '  +3.09%  '
$ws.Range('E36').Value = '  +15.80%  '
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.12'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.04'
$ws.Range('E40').Value = '  +4.80%  '
$ws.Range('E41').Value = '  -3.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '133.80'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.290'
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.85'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.25'
$ws.Range('E47').Value = '  -5.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.21'
$ws.Range('E48').Value = '  +3.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '21.09'
$ws.Range('E49').Value = '  -5.30%  '
$ws.Range('D50').Value = '2.136.01'
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = '3.645.80'
$ws.Range('E51').Value = '  +2.00%  '
